# Apply updated dSF (column F) values — repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -5
$ws.Range("F3").Value  = 7
$ws.Range("F5").Value  = -7
$ws.Range("F6").Value  = -1
$ws.Range("F8").Value  = -5
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 9
$ws.Range("F18").Value = -3
$ws.Range("F23").Value = 2
